# Re-apply per-leve market-price / profit recalculation figures.
# Source: scheduled market-data refresh (see commit message).
# Workbook has no formulas -- every Hxx..Nxx figure below is a plain
# literal pulled from the latest pricing snapshot, grouped by sheet/row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
# row 12
$ws.Range("H12").Value = 501
$ws.Range("J12").Value = 501
$ws.Range("L12").Value = 501
$ws.Range("N12").Value = -841
# row 18
$ws.Range("H18").Value = 17333.334
$ws.Range("I18").Value = 17333.334
$ws.Range("K18").Value = 17333.334
$ws.Range("M18").Value = -17049.334
# row 32
$ws.Range("H32").Value = 2122.5386
$ws.Range("I32").Value = 2000
$ws.Range("K32").Value = 2000
$ws.Range("M32").Value = -1674
# row 46
$ws.Range("H46").Value = 20650
$ws.Range("J46").Value = 29300
$ws.Range("L46").Value = 87900
$ws.Range("N46").Value = -88138
# row 51
$ws.Range("H51").Value = 5625
$ws.Range("J51").Value = 5625
$ws.Range("L51").Value = 5625
$ws.Range("N51").Value = -6593
# row 60
$ws.Range("H60").Value = 20650
$ws.Range("J60").Value = 29300
$ws.Range("L60").Value = 87900
$ws.Range("N60").Value = -88868
# row 62
$ws.Range("H62").Value = 4995
$ws.Range("I62").Value = 4990
$ws.Range("K62").Value = 4990
$ws.Range("M62").Value = -4366
# row 65
$ws.Range("H65").Value = 4995
$ws.Range("I65").Value = 4990
$ws.Range("K65").Value = 24950
$ws.Range("M65").Value = -21830
# row 98
$ws.Range("H98").Value = 62500428
$ws.Range("I98").Value = 62500428
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 62500428
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -62498930
$ws.Range("N98").Value = ""
# row 113
$ws.Range("H113").Value = 1499.6666
$ws.Range("I113").Value = 1499.5
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1499.5
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1754.5
$ws.Range("N113").Value = -8008
# row 122
$ws.Range("H122").Value = 62500428
$ws.Range("I122").Value = 62500428
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 187501284
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -187498834
$ws.Range("N122").Value = ""
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 924.7143
$ws.Range("I2").Value = 962.1667
$ws.Range("K2").Value = 962.1667
$ws.Range("M2").Value = -849.1667
# row 30
$ws.Range("H30").Value = 5417.4165
$ws.Range("I30").Value = 5417.4165
$ws.Range("K30").Value = 5417.4165
$ws.Range("M30").Value = -5267.4165
# row 50
$ws.Range("H50").Value = 7400.3335
$ws.Range("I50").Value = 1074
$ws.Range("K50").Value = 1074
$ws.Range("M50").Value = -360
# row 110
$ws.Range("H110").Value = 1185.875
$ws.Range("I110").Value = 1420
$ws.Range("K110").Value = 1420
$ws.Range("M110").Value = 625
# row 116
$ws.Range("H116").Value = 924.7143
$ws.Range("I116").Value = 962.1667
$ws.Range("K116").Value = 962.1667
$ws.Range("M116").Value = 1331.8333
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 924.7143
$ws.Range("I3").Value = 962.1667
$ws.Range("K3").Value = 962.1667
$ws.Range("M3").Value = -848.1667
# row 80
$ws.Range("H80").Value = 1235.25
$ws.Range("I80").Value = 982.3333
$ws.Range("K80").Value = 982.3333
$ws.Range("M80").Value = 15.66669999999999
# row 83
$ws.Range("H83").Value = 1235.25
$ws.Range("I83").Value = 982.3333
$ws.Range("K83").Value = 4911.6665
$ws.Range("M83").Value = 80.33349999999973
$ws = $wb.Worksheets.Item("CRP")
# row 5
$ws.Range("H5").Value = 216.08333
$ws.Range("I5").Value = 217.54546
$ws.Range("K5").Value = 217.54546
$ws.Range("M5").Value = -105.54546
# row 16
$ws.Range("H16").Value = 4999.25
$ws.Range("I16").Value = 4999.3335
$ws.Range("J16").Value = 4999
$ws.Range("K16").Value = 4999.3335
$ws.Range("L16").Value = 4999
$ws.Range("M16").Value = -4712.3335
$ws.Range("N16").Value = -5573
# row 22
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = -1700
# row 74
$ws.Range("H74").Value = 29499.5
$ws.Range("J74").Value = 29499.5
$ws.Range("L74").Value = 29499.5
$ws.Range("N74").Value = -31247.5
# row 77
$ws.Range("H77").Value = 29499.5
$ws.Range("J77").Value = 29499.5
$ws.Range("L77").Value = 88498.5
$ws.Range("N77").Value = -97234.5
# row 113
$ws.Range("H113").Value = 4999.25
$ws.Range("I113").Value = 4999.3335
$ws.Range("J113").Value = 4999
$ws.Range("K113").Value = 4999.3335
$ws.Range("L113").Value = 4999
$ws.Range("M113").Value = -2829.3335
$ws.Range("N113").Value = -9339
$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 17.9
$ws.Range("I2").Value = 17.375
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 104.25
$ws.Range("L2").Value = 120
$ws.Range("M2").Value = 8.75
$ws.Range("N2").Value = -346
# row 6
$ws.Range("H6").Value = 111162.445
$ws.Range("J6").Value = 333397.34
$ws.Range("L6").Value = 1000192.02
$ws.Range("N6").Value = -1000418.02
# row 7
$ws.Range("H7").Value = 157.8
$ws.Range("I7").Value = 227.2
$ws.Range("K7").Value = 681.5999999999999
$ws.Range("M7").Value = -569.5999999999999
# row 12
$ws.Range("H12").Value = 5.857143
$ws.Range("J12").Value = 5.857143
$ws.Range("L12").Value = 17.571429
$ws.Range("N12").Value = -363.571429
# row 17
$ws.Range("H17").Value = 563
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2831
# row 34
$ws.Range("H34").Value = 689.7
$ws.Range("J34").Value = 2096.3333
$ws.Range("L34").Value = 6288.999899999999
$ws.Range("N34").Value = -6456.999899999999
# row 39
$ws.Range("H39").Value = 6775
$ws.Range("I39").Value = 5033.3335
$ws.Range("J39").Value = 12000
$ws.Range("K39").Value = 15100.0005
$ws.Range("L39").Value = 36000
$ws.Range("M39").Value = -14806.0005
$ws.Range("N39").Value = -36588
# row 40
$ws.Range("H40").Value = 250
$ws.Range("J40").Value = 250
$ws.Range("L40").Value = 1000
$ws.Range("N40").Value = -1138
# row 46
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = ""
# row 55
$ws.Range("H55").Value = 9233.333000000001
$ws.Range("I55").Value = 1200
$ws.Range("J55").Value = 13250
$ws.Range("K55").Value = 3600
$ws.Range("L55").Value = 39750
$ws.Range("M55").Value = -3423
$ws.Range("N55").Value = -40104
# row 122
$ws.Range("H122").Value = 874.25
$ws.Range("I122").Value = 499
$ws.Range("J122").Value = 999.3333
$ws.Range("K122").Value = 4491
$ws.Range("L122").Value = 8993.9997
$ws.Range("M122").Value = -2041
$ws.Range("N122").Value = -13893.9997
$ws = $wb.Worksheets.Item("GSM")
# row 95
$ws.Range("H95").Value = 28975
$ws.Range("J95").Value = 28975
$ws.Range("L95").Value = 28975
$ws.Range("N95").Value = -34467
# row 122
$ws.Range("H122").Value = 2837.4285
$ws.Range("I122").Value = 2868
$ws.Range("K122").Value = 8604
$ws.Range("M122").Value = -6154
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 2066.6667
$ws.Range("I22").Value = 1850
$ws.Range("K22").Value = 1850
$ws.Range("M22").Value = -1555
# row 27
$ws.Range("H27").Value = 2066.6667
$ws.Range("I27").Value = 1850
$ws.Range("K27").Value = 1850
$ws.Range("M27").Value = -1743
# row 40
$ws.Range("H40").Value = 4396.8237
$ws.Range("I40").Value = 3730.75
$ws.Range("J40").Value = 4988.8887
$ws.Range("K40").Value = 3730.75
$ws.Range("L40").Value = 4988.8887
$ws.Range("M40").Value = -3594.75
$ws.Range("N40").Value = -5260.8887
# row 41
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = ""
$ws = $wb.Worksheets.Item("WVR")
# row 44
$ws.Range("H44").Value = 24990
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").Value = ""
# row 126
$ws.Range("H126").Value = 8375.75
$ws.Range("I126").Value = 8375.75
$ws.Range("K126").Value = 25127.25
$ws.Range("M126").Value = -22657.25
# row 132
$ws.Range("H132").Value = 4260.375
$ws.Range("I132").Value = 4260.375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12781.125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10251.125
$ws.Range("N132").Value = ""
# row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = ""
